# Automated daily data refresh: append 08-sep to "Prix Spot" and the
# 2025-09-06 / 2025-09-07 rows to "Gaz" and "CO2".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": append column CI ("08-sep") after CH ("07-sep").
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Clone the header cell's look (bold font + border + centered) from the
# previous day's column, then overwrite with the new day's label.
$wsPrix.Cells.Item(1, 86).Copy($wsPrix.Cells.Item(1, 87))
$wsPrix.Cells.Item(1, 87).Value = "08-sep"

$prixValues = @(
    66.76000000000001,
    63.5,
    63.39,
    37.78,
    31.86,
    55.7,
    77.98,
    69.61,
    108.36,
    98.20999999999999,
    85.83,
    63,
    52.01,
    38.56,
    37.21,
    27.34,
    38.72,
    48.64,
    73.08,
    91.15000000000001,
    98.20999999999999,
    95.26000000000001,
    89.34999999999999,
    73.08
)

for ($i = 0; $i -lt $prixValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 87).Value = $prixValues[$i]
}

# ---------------------------------------------------------------------
# Helper: write a literal (non-date-parsed) text value into a target
# cell without leaving any explicit style on it. Entering an ISO-style
# date string straight into .Value gets auto-detected and converted to
# a date serial, so instead we compute it as a throwaway string formula
# in a scratch cell, copy the *result* over (which lands as plain text
# with no number format), then clear the scratch cell again.
# ---------------------------------------------------------------------
function Set-LiteralText($scratchCell, $targetCell, $text) {
    $scratchCell.Formula = '="' + $text + '"'
    $scratchCell.Copy($targetCell)
    $scratchCell.Clear()
}

# ---------------------------------------------------------------------
# Sheet "Gaz": append rows 84 (2025-09-06) and 85 (2025-09-07).
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$gazScratch = $wsGaz.Cells.Item(500, 500)

Set-LiteralText $gazScratch $wsGaz.Cells.Item(84, 1) "2025-09-06"
$wsGaz.Cells.Item(84, 2).Value = 31

Set-LiteralText $gazScratch $wsGaz.Cells.Item(85, 1) "2025-09-07"
$wsGaz.Cells.Item(85, 2).Value = 31

# ---------------------------------------------------------------------
# Sheet "CO2": append rows 84 (2025-09-06) and 85 (2025-09-07).
# ---------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")
$co2Scratch = $wsCO2.Cells.Item(500, 500)

Set-LiteralText $co2Scratch $wsCO2.Cells.Item(84, 1) "2025-09-06"
$wsCO2.Cells.Item(84, 2).Value = 75.59

Set-LiteralText $co2Scratch $wsCO2.Cells.Item(85, 1) "2025-09-07"
$wsCO2.Cells.Item(85, 2).Value = 75.59
